$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1785.6
$ws.Range("I28").Value = 2299.8333
$ws.Range("K28").Value = 2299.8333
$ws.Range("M28").Value = -1814.8333
$ws.Range("H51").Value = 15847.625
$ws.Range("J51").Value = 9540.286
$ws.Range("L51").Value = 9540.286
$ws.Range("N51").Value = -10508.286
$ws.Range("H88").Value = 15918022
$ws.Range("I88").Value = 37040610
$ws.Range("J88").Value = 76081
$ws.Range("K88").Value = 37040610
$ws.Range("L88").Value = 76081
$ws.Range("M88").Value = -37040204
$ws.Range("N88").Value = -76893
$ws.Range("H91").Value = 15918022
$ws.Range("I91").Value = 37040610
$ws.Range("J91").Value = 76081
$ws.Range("K91").Value = 37040610
$ws.Range("L91").Value = 76081
$ws.Range("M91").Value = -37039206
$ws.Range("N91").Value = -78889
$ws.Range("H101").Value = 741.5714
$ws.Range("J101").Value = 2685
$ws.Range("L101").Value = 8055
$ws.Range("N101").Value = -11299
$ws.Range("H113").Value = 33339816
$ws.Range("I113").Value = 2220.6667
$ws.Range("K113").Value = 2220.6667
$ws.Range("M113").Value = 1033.3333
$ws.Range("H132").Value = 1484.2778
$ws.Range("I132").Value = 1469.5714
$ws.Range("K132").Value = 4408.7142
$ws.Range("M132").Value = -1878.7142
$ws.Range("H137").Value = 4005.2222
$ws.Range("I137").Value = 6320.2
$ws.Range("J137").Value = 3114.8462
$ws.Range("K137").Value = 18960.6
$ws.Range("L137").Value = 9344.5386
$ws.Range("M137").Value = -16410.6
$ws.Range("N137").Value = -14444.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1546880.4
$ws.Range("I32").Value = 1648353.5
$ws.Range("K32").Value = 1648353.5
$ws.Range("M32").Value = -1648066.5
$ws.Range("H45").Value = 5625.2856
$ws.Range("I45").Value = 2601.8333
$ws.Range("K45").Value = 2601.8333
$ws.Range("M45").Value = -2224.8333
$ws.Range("H61").Value = 5020.2393
$ws.Range("I61").Value = 2058.2942
$ws.Range("K61").Value = 2058.2942
$ws.Range("M61").Value = -1846.2942
$ws.Range("H63").Value = 2408.2856
$ws.Range("I63").Value = 2369.5
$ws.Range("J63").Value = 2460
$ws.Range("K63").Value = 2369.5
$ws.Range("L63").Value = 2460
$ws.Range("M63").Value = -1683.5
$ws.Range("N63").Value = -3832
$ws.Range("H66").Value = 2408.2856
$ws.Range("I66").Value = 2369.5
$ws.Range("J66").Value = 2460
$ws.Range("K66").Value = 11847.5
$ws.Range("L66").Value = 12300
$ws.Range("M66").Value = -8415.5
$ws.Range("N66").Value = -19164
$ws.Range("H136").Value = 5020.2393
$ws.Range("I136").Value = 2058.2942
$ws.Range("K136").Value = 6174.882599999999
$ws.Range("M136").Value = -3624.882599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 66667508
$ws.Range("I64").Value = 111111780
$ws.Range("J64").Value = 1100
$ws.Range("K64").Value = 111111780
$ws.Range("L64").Value = 1100
$ws.Range("M64").Value = -111111555
$ws.Range("N64").Value = -1550
$ws.Range("H67").Value = 66667508
$ws.Range("I67").Value = 111111780
$ws.Range("J67").Value = 1100
$ws.Range("K67").Value = 111111780
$ws.Range("L67").Value = 1100
$ws.Range("M67").Value = -111111000
$ws.Range("N67").Value = -2660
$ws.Range("H94").Value = 1425.7142
$ws.Range("I94").Value = 582.53845
$ws.Range("J94").Value = 2795.875
$ws.Range("K94").Value = 582.53845
$ws.Range("L94").Value = 2795.875
$ws.Range("M94").Value = -131.53845
$ws.Range("N94").Value = -3697.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5976.3
$ws.Range("J62").Value = 4624
$ws.Range("L62").Value = 4624
$ws.Range("N62").Value = -5872
$ws.Range("H65").Value = 5976.3
$ws.Range("J65").Value = 4624
$ws.Range("L65").Value = 23120
$ws.Range("N65").Value = -29360
$ws.Range("H86").Value = 8935629
$ws.Range("I86").Value = 20839000
$ws.Range("K86").Value = 20839000
$ws.Range("M86").Value = -20837877
$ws.Range("H89").Value = 8935629
$ws.Range("I89").Value = 20839000
$ws.Range("K89").Value = 104195000
$ws.Range("M89").Value = -104189384
$ws.Range("H99").Value = 6796.8
$ws.Range("I99").Value = 7366.5835
$ws.Range("K99").Value = 7366.5835
$ws.Range("M99").Value = -5868.5835
$ws.Range("H126").Value = 6796.8
$ws.Range("I126").Value = 7366.5835
$ws.Range("K126").Value = 22099.7505
$ws.Range("M126").Value = -19629.7505
$ws.Range("H132").Value = 9528529
$ws.Range("I132").Value = 1695.9048
$ws.Range("J132").Value = 19055362
$ws.Range("K132").Value = 5087.7144
$ws.Range("L132").Value = 57166086
$ws.Range("M132").Value = -2557.7144
$ws.Range("N132").Value = -57171146

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 9000
$ws.Range("M76").Value = -8617
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 9000
$ws.Range("M79").Value = -7674
$ws.Range("H107").Value = 13333606
$ws.Range("J107").Value = 16666908
$ws.Range("L107").Value = 50000724
$ws.Range("N107").Value = -50004564

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 667
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 667
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H61").Value = 4830.0713
$ws.Range("I61").Value = 3010.7334
$ws.Range("J61").Value = 6929.3076
$ws.Range("K61").Value = 3010.7334
$ws.Range("L61").Value = 6929.3076
$ws.Range("M61").Value = -2808.7334
$ws.Range("N61").Value = -7333.3076
$ws.Range("H113").Value = 4830.0713
$ws.Range("I113").Value = 3010.7334
$ws.Range("J113").Value = 6929.3076
$ws.Range("K113").Value = 3010.7334
$ws.Range("L113").Value = 6929.3076
$ws.Range("M113").Value = -840.7334000000001
$ws.Range("N113").Value = -11269.3076
$ws.Range("H132").Value = 9265127
$ws.Range("I132").Value = 19232828
$ws.Range("J132").Value = 9404.679
$ws.Range("K132").Value = 57698484
$ws.Range("L132").Value = 28214.037
$ws.Range("M132").Value = -57695954
$ws.Range("N132").Value = -33274.037
$ws.Range("H136").Value = 10739.92
$ws.Range("I136").Value = 2249.6667
$ws.Range("J136").Value = 13421.053
$ws.Range("K136").Value = 6749.000100000001
$ws.Range("L136").Value = 40263.159
$ws.Range("M136").Value = -4199.000100000001
$ws.Range("N136").Value = -45363.159

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14292129
$ws.Range("I81").Value = 2523.7
$ws.Range("J81").Value = 50016144
$ws.Range("K81").Value = 5047.4
$ws.Range("L81").Value = 100032288
$ws.Range("M81").Value = -3986.4
$ws.Range("N81").Value = -100034410
$ws.Range("H84").Value = 14292129
$ws.Range("I84").Value = 2523.7
$ws.Range("J84").Value = 50016144
$ws.Range("K84").Value = 25237
$ws.Range("L84").Value = 500161440
$ws.Range("M84").Value = -19933
$ws.Range("N84").Value = -500172048
$ws.Range("H100").Value = 885.36
$ws.Range("J100").Value = 1015.9167
$ws.Range("L100").Value = 2031.8334
$ws.Range("N100").Value = -3113.8334
$ws.Range("H136").Value = 21765158
$ws.Range("I136").Value = 40000916
$ws.Range("J136").Value = 55925.19
$ws.Range("K136").Value = 120002748
$ws.Range("L136").Value = 167775.57
$ws.Range("M136").Value = -120000198
$ws.Range("N136").Value = -172875.57
